$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'68.322.01"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.75%  "

# Row 3
$ws.Range("D3").Value = "'3.925.20"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.42%  "

# Row 4
$ws.Range("E4").Value = "  +0.14%  "

# Row 5
$ws.Range("D5").Value = "'487.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.38%  "

# Row 6
$ws.Range("D6").Value = "'147.81"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.22%  "

# Row 7
$ws.Range("E7").Value = "  -0.89%  "

# Row 8
$ws.Range("E8").Value = "  -0.03%  "

# Row 9
$ws.Range("D9").Value = "'0.733"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.03%  "

# Row 10
$ws.Range("E10").Value = "  +1.68%  "

# Row 11
$ws.Range("E11").Value = "  +3.96%  "

# Row 12
$ws.Range("D12").Value = "'43.11"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.51%  "

# Row 13
$ws.Range("D13").Value = "'10.74"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.82%  "

# Row 14
$ws.Range("D14").Value = "'4.554.94"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.33%  "

# Row 15
$ws.Range("D15").Value = "'3.920.77"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.36%  "

# Row 16
$ws.Range("D16").Value = "'14.44"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.43%  "

# Row 17
$ws.Range("E17").Value = "  -0.68%  "

# Row 19
$ws.Range("D19").Value = "'1.13"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.09%  "

# Row 20
$ws.Range("D20").Value = "'68.416.32"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.36%  "

# Row 21
$ws.Range("D21").Value = "'441.98"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.83%  "

# Row 22
$ws.Range("E22").Value = "  +2.89%  "

# Row 23
$ws.Range("D23").Value = "'15.16"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.36%  "

# Row 24
$ws.Range("D24").Value = "'88.22"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.61%  "

# Row 25
$ws.Range("D25").Value = "'11.44"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +17.17%  "

# Row 26
$ws.Range("D26").Value = "'11.52"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +13.17%  "

# Row 27
$ws.Range("D27").Value = "'3.64"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.49%  "

# Row 28
$ws.Range("D28").Value = "'38.68"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.75%  "

# Row 29
$ws.Range("E29").Value = "  -1.10%  "

# Row 30
$ws.Range("D30").Value = "'717.82"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.72%  "

# Row 31
$ws.Range("D31").Value = "'13.83"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.04%  "

# Row 32
$ws.Range("E32").Value = "  -0.90%  "

# Row 33
$ws.Range("E33").Value = "  +2.88%  "

# Row 34
$ws.Range("D34").Value = "'6.35"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +18.43%  "

# Row 35
$ws.Range("D35").Value = "'42.29"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.95%  "

# Row 36
$ws.Range("D36").Value = "'0.0₃0877"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +11.65%  "

# Row 37
$ws.Range("D37").Value = "'61.47"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +6.25%  "

# Row 38
$ws.Range("D38").Value = "'0.421"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +25.14%  "

# Row 39
$ws.Range("E39").Value = "  -2.20%  "

# Row 40
$ws.Range("E40").Value = "  +18.86%  "

# Row 41
$ws.Range("B41").Value = "ThetaToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D41").Value = "'3.33"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +8.64%  "

# Row 42
$ws.Range("B42").Value = "Dai"
$ws.Range("C42").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D42").Value = "'0.998"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.08%  "

# Row 43
$ws.Range("D43").Value = "'0.0481"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.77%  "

# Row 44
$ws.Range("D44").Value = "'2.93"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.93%  "

# Row 45
$ws.Range("E45").Value = "  -0.07%  "

# Row 46
$ws.Range("E46").Value = "  +0.01%  "

# Row 47
$ws.Range("D47").Value = "'3.32"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.93%  "

# Row 48
$ws.Range("E48").Value = "  -1.80%  "

# Row 49
$ws.Range("D49").Value = "'0.0₆0351"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +32.03%  "

# Row 50
$ws.Range("D50").Value = "'2.14"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.20%  "

# Row 51
$ws.Range("D51").Value = "'146.15"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.33%  "
